$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "57.805.87"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.83%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.448.21"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.40%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "509.95"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -2.17%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "129.56"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -1.93%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.549"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -1.89%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.468.54"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.89%  "
$ws.Range("B10").Value = "TRON"
$ws.Range("C10").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.156"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.23%  "
$ws.Range("B11").Value = "Dogecoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0960"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -3.76%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -3.57%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.327"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -5.01%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.885.14"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.30%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "57.746.98"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.79%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.84"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -2.82%  "
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -2.72%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.462.68"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.92%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.50"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -3.58%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "317.64"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.80%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.11"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -1.93%  "
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.15%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.92"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +2.37%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "62.94"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -2.21%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.399"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -3.10%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.990"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.76%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.37%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.22"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -2.16%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "168.71"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +1.29%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0728"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -4.13%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -2.83%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -1.71%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -3.23%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.09%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.996"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.21%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.75"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -2.11%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -4.86%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -2.92%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.63"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.10%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -2.61%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.761"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -4.30%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "270.64"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -2.44%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.00"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -1.61%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -3.81%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.583"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -2.52%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.29%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "120.10"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -5.90%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0486"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.91%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "17.26"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -4.22%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -2.42%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "16.64"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -3.33%  "
